# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 13
$ws1.Range("F6").Value  = 544
$ws1.Range("F7").Value  = 1677
$ws1.Range("F8").Value  = 20
$ws1.Range("F10").Value = 28
$ws1.Range("F11").Value = 1630
$ws1.Range("F13").Value = 68
$ws1.Range("F14").Value = 402
$ws1.Range("F15").Value = 264
$ws1.Range("F18").Value = 23
$ws1.Range("F19").Value = 31
$ws1.Range("F20").Value = 50
$ws1.Range("F21").Value = 223
$ws1.Range("F22").Value = 291
$ws1.Range("F23").Value = 159
$ws1.Range("F24").Value = 224
$ws1.Range("F25").Value = 226

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 13
$ws4.Range("F6").Value  = 544
$ws4.Range("F7").Value  = 1677
$ws4.Range("F9").Value  = 20
$ws4.Range("F11").Value = 28
$ws4.Range("F12").Value = 1630
$ws4.Range("F14").Value = 68
$ws4.Range("F15").Value = 402
$ws4.Range("F16").Value = 264
$ws4.Range("F19").Value = 23
$ws4.Range("F20").Value = 31
$ws4.Range("F21").Value = 50
$ws4.Range("F22").Value = 223
$ws4.Range("F23").Value = 291
$ws4.Range("F24").Value = 159
$ws4.Range("F25").Value = 224
$ws4.Range("F26").Value = 226
